$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

$ws.Range("B2").Value = "Tue Jan 28 21:48:45 EST 2025"
$ws.Range("B3").Value = "Tue Jan 28 21:48:55 EST 2025"
$ws.Range("B4").Value = "Tue Jan 28 21:49:04 EST 2025"
$ws.Range("B5").Value = "Tue Jan 28 21:49:13 EST 2025"
$ws.Range("B6").Value = "Tue Jan 28 21:49:23 EST 2025"
$ws.Range("B7").Value = "Tue Jan 28 21:49:33 EST 2025"
$ws.Range("B8").Value = "Tue Jan 28 21:49:42 EST 2025"
$ws.Range("B9").Value = "Tue Jan 28 21:49:51 EST 2025"
$ws.Range("B10").Value = "Tue Jan 28 21:50:01 EST 2025"
$ws.Range("B11").Value = "Tue Jan 28 21:50:11 EST 2025"
$ws.Range("B12").Value = "Tue Jan 28 21:50:20 EST 2025"
$ws.Range("B13").Value = "Tue Jan 28 21:50:30 EST 2025"
$ws.Range("B14").Value = "Tue Jan 28 21:50:39 EST 2025"
$ws.Range("B15").Value = "Tue Jan 28 21:50:48 EST 2025"
$ws.Range("B16").Value = "Tue Jan 28 21:50:58 EST 2025"
$ws.Range("B17").Value = "Tue Jan 28 21:51:07 EST 2025"
$ws.Range("B18").Value = "Tue Jan 28 21:51:16 EST 2025"
$ws.Range("B19").Value = "Tue Jan 28 21:51:26 EST 2025"
$ws.Range("B20").Value = "Tue Jan 28 21:51:35 EST 2025"
$ws.Range("B21").Value = "Tue Jan 28 21:51:44 EST 2025"
$ws.Range("B22").Value = "Tue Jan 28 21:51:54 EST 2025"
$ws.Range("B23").Value = "Tue Jan 28 21:52:03 EST 2025"
$ws.Range("B24").Value = "Tue Jan 28 21:52:12 EST 2025"
$ws.Range("B25").Value = "Tue Jan 28 21:52:22 EST 2025"
$ws.Range("B26").Value = "Tue Jan 28 21:52:31 EST 2025"
$ws.Range("B27").Value = "Tue Jan 28 21:52:40 EST 2025"
$ws.Range("B28").Value = "Tue Jan 28 21:52:50 EST 2025"
$ws.Range("B29").Value = "Tue Jan 28 21:52:59 EST 2025"
$ws.Range("B30").Value = "Tue Jan 28 21:53:08 EST 2025"
$ws.Range("B31").Value = "Tue Jan 28 21:53:18 EST 2025"
$ws.Range("B32").Value = "Tue Jan 28 21:53:27 EST 2025"
$ws.Range("B33").Value = "Tue Jan 28 21:53:36 EST 2025"
$ws.Range("B34").Value = "Tue Jan 28 21:53:46 EST 2025"
$ws.Range("B35").Value = "Tue Jan 28 21:53:55 EST 2025"
$ws.Range("B36").Value = "Tue Jan 28 21:54:04 EST 2025"
$ws.Range("B37").Value = "Tue Jan 28 21:54:14 EST 2025"
$ws.Range("B38").Value = "Tue Jan 28 21:54:24 EST 2025"
$ws.Range("B39").Value = "Tue Jan 28 21:54:33 EST 2025"
$ws.Range("B40").Value = "Tue Jan 28 21:54:42 EST 2025"
$ws.Range("B44").Value = "Tue Jan 28 21:54:52 EST 2025"
$ws.Range("B45").Value = "Tue Jan 28 21:55:01 EST 2025"
$ws.Range("B46").Value = "Tue Jan 28 21:55:10 EST 2025"
$ws.Range("B47").Value = "Tue Jan 28 21:55:20 EST 2025"
$ws.Range("B48").Value = "Tue Jan 28 21:55:30 EST 2025"
$ws.Range("B49").Value = "Tue Jan 28 21:55:40 EST 2025"
$ws.Range("B50").Value = "Tue Jan 28 21:55:49 EST 2025"
$ws.Range("B51").Value = "Tue Jan 28 21:55:59 EST 2025"
$ws.Range("B52").Value = "Tue Jan 28 21:56:08 EST 2025"
$ws.Range("B53").Value = "Tue Jan 28 21:56:18 EST 2025"
$ws.Range("B54").Value = "Tue Jan 28 21:56:27 EST 2025"
$ws.Range("B55").Value = "Tue Jan 28 21:56:37 EST 2025"
$ws.Range("B56").Value = "Tue Jan 28 21:56:46 EST 2025"
$ws.Range("B57").Value = "Tue Jan 28 21:56:55 EST 2025"
$ws.Range("B58").Value = "Tue Jan 28 21:57:05 EST 2025"
$ws.Range("B59").Value = "Tue Jan 28 21:57:14 EST 2025"
$ws.Range("B60").Value = "Tue Jan 28 21:57:24 EST 2025"
$ws.Range("B61").Value = "Tue Jan 28 21:57:33 EST 2025"
